$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3) ---
$ws.Range("A3").Value = "Designator"
$ws.Range("B3").Value = "Quantity"
$ws.Range("C3").Value = "Description"
$ws.Range("D3").Value = "Part Number"
$ws.Range("E3").Value = "DigiKey Part Number"
$ws.Range("F3").Value = "Package"

# --- Row 4: MCU ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "IC MCU 32BIT 128KB FLASH 48QFPN"
$ws.Range("D4").Value = "STM32F102CBU6"
$ws.Range("E4").Value = "497-17381-ND"
$ws.Range("F4").Value = "UFQFPN 48"

# --- Row 5: USB connector ---
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "USB connector"

# --- Row 6: debug header connector ---
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "CONN HEADER SMD 14POS 1.27MM"
$ws.Range("D6").Value = "20021121-00014T4LF"
$ws.Range("E6").Value = "609-3730-ND"

# --- Row 7: 3.3V LDO ---
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "3.3V LDO"

# --- Row 8: Red LED ---
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "Red LED"

# --- Row 9: Green LED ---
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "Green LED"

# --- Row 10: Slide switch (replaces "Target Power Switch") ---
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Slide Switch DPDT Surface Mount, Right Angle"
$ws.Range("D10").Value = "JS202011JAQN"
$ws.Range("E10").Value = "CKN10722CT-ND"

# --- Column widths (target stored widths: C=32.140625, D:E=23.7109375, F=10.7109375;
#     the engine snaps ColumnWidth assignments to 1/6-character pixel grid, so these
#     inputs are chosen to land on the closest achievable snap point) ---
$ws.Columns.Item(3).ColumnWidth = 31.33
$ws.Columns.Item(4).ColumnWidth = 22.83
$ws.Columns.Item(5).ColumnWidth = 22.83
$ws.Columns.Item(6).ColumnWidth = 9.83

# --- Selection ---
$null = $ws.Range("N5").Select()
